$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Numara: -> student number
$ws.Range("J5").Value = 20215070019

# Ad Soyad: -> full name
$ws.Range("J6").Value = "KÜBRA ÇABUK"

# Bölüm: -> department
$ws.Range("J7").Value = "YBS"

# Concatenation formula demo cell
$ws.Range("D10").Formula = "=B4&F4"

# Update selection to match the authored state
$ws.Range("J7:L7").Select()
